$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 24 (1-indexed): "Ratio" row in the RETENTION section -> value 0.8
$cell = $t.Rows.Item(24).Cells.Item(2)
$cell.Range.Text = "0.8"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 12
$cell.Range.Font.SizeBi = 12

# Row 44 (1-indexed): "Answer Recall Lenient (ARL)" -> value 0.4285
$cell = $t.Rows.Item(44).Cells.Item(2)
$cell.Range.Text = "0.4285"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 12
$cell.Range.Font.SizeBi = 12

# Row 45 (1-indexed): "Answer Recall Strict (ARS)" -> value 0.1428
$cell = $t.Rows.Item(45).Cells.Item(2)
$cell.Range.Text = "0.1428"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 12
$cell.Range.Font.SizeBi = 12

# Row 46 (1-indexed): "Answer Recall Average (ARA)" -> value 0.2856
$cell = $t.Rows.Item(46).Cells.Item(2)
$cell.Range.Text = "0.2856"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 12
$cell.Range.Font.SizeBi = 12
